# Auto-generated Excel COM-interop script
# Applies updated crypto price/volume values per commit:
# "Updated cryptos list on Wed Aug 16 11:15:17 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "29.248.68"
Set-TextCell "E2" "  -0.44%  "
Set-TextCell "D3" "1.829.89"
Set-TextCell "D4" "1.002"
Set-TextCell "E4" "  +0.23%  "
Set-TextCell "D5" "235.28"
Set-TextCell "E5" "  -1.62%  "
Set-TextCell "D6" "0.6031"
Set-TextCell "E6" "  -3.77%  "
Set-TextCell "D7" "1.004"
Set-TextCell "E7" "  +0.33%  "
Set-TextCell "D8" "0.07052"
Set-TextCell "E8" "  -5.12%  "
Set-TextCell "D9" "0.2798"
Set-TextCell "E9" "  -3.24%  "
Set-TextCell "D10" "23.56"
Set-TextCell "E10" "  -5.51%  "
Set-TextCell "D11" "0.07650"
Set-TextCell "E11" "  -0.93%  "
Set-TextCell "D12" "1.833.16"
Set-TextCell "E12" "  -0.59%  "
Set-TextCell "D13" "4.794"
Set-TextCell "E13" "  -3.65%  "
Set-TextCell "D14" "0.6291"
Set-TextCell "E14" "  -6.63%  "
Set-TextCell "D15" "0.000009877"
Set-TextCell "E15" "  -3.93%  "
Set-TextCell "D16" "2.084.26"
Set-TextCell "E16" "  -0.31%  "
Set-TextCell "D17" "79.08"
Set-TextCell "E17" "  -3.29%  "
Set-TextCell "B18" "WrappedBTC"
Set-TextCell "C18" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D18" "29.259.52"
Set-TextCell "E18" "  -0.57%  "
Set-TextCell "B19" "Uniswap"
Set-TextCell "C19" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D19" "5.837"
Set-TextCell "E19" "  -5.95%  "
Set-TextCell "D20" "224.20"
Set-TextCell "E20" "  -4.29%  "
Set-TextCell "E21" "  +0.22%  "
Set-TextCell "E22" "  -4.95%  "
Set-TextCell "D23" "7.004"
Set-TextCell "E23" "  -3.95%  "
Set-TextCell "E24" "  +0.26%  "
Set-TextCell "D25" "156.42"
Set-TextCell "E25" "  -0.83%  "
Set-TextCell "D26" "0.1303"
Set-TextCell "E26" "  -3.02%  "
Set-TextCell "D27" "7.980"
Set-TextCell "E27" "  -6.15%  "
Set-TextCell "D28" "16.64"
Set-TextCell "E28" "  -3.77%  "
Set-TextCell "D29" "1.482"
Set-TextCell "E29" "  +1.05%  "
Set-TextCell "D30" "0.06517"
Set-TextCell "E30" "  -9.94%  "
Set-TextCell "E31" "  -1.95%  "
Set-TextCell "D32" "3.837"
Set-TextCell "E32" "  -4.78%  "
Set-TextCell "D33" "3.796"
Set-TextCell "E33" "  -5.92%  "
Set-TextCell "D34" "1.108"
Set-TextCell "E34" "  -2.72%  "
Set-TextCell "D35" "1.730"
Set-TextCell "E35" "  -4.74%  "
Set-TextCell "D36" "0.6463"
Set-TextCell "E36" "  -7.34%  "
Set-TextCell "D37" "2.543"
Set-TextCell "E37" "  -1.33%  "
Set-TextCell "D38" "1.215.75"
Set-TextCell "E38" "  -1.55%  "
Set-TextCell "D39" "2.741"
Set-TextCell "E39" "  -2.71%  "
Set-TextCell "D40" "0.01747"
Set-TextCell "E40" "  -5.23%  "
Set-TextCell "D41" "6.564"
Set-TextCell "E41" "  -5.32%  "
Set-TextCell "D42" "0.8966"
Set-TextCell "E42" "  -6.63%  "
Set-TextCell "E43" "  +0.26%  "
Set-TextCell "D44" "1.998.68"
Set-TextCell "E44" "  -0.02%  "
Set-TextCell "D45" "100.35"
Set-TextCell "E45" "  -0.48%  "
Set-TextCell "D46" "62.64"
Set-TextCell "E46" "  -4.17%  "
Set-TextCell "D47" "0.00000000116"
Set-TextCell "E47" "  -3.52%  "
Set-TextCell "D48" "8.560"
Set-TextCell "E48" "  -3.66%  "
Set-TextCell "D49" "1.579"
Set-TextCell "E49" "  -8.09%  "
Set-TextCell "D50" "0.4549"
Set-TextCell "E50" "  -0.58%  "
Set-TextCell "D51" "0.05500"
Set-TextCell "E51" "  -2.80%  "
